$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A labels for rows 7-21 to their new values (rows 7-10 just swap,
# rows 11-21 shift up by one row with row 11 replaced by "MATTEr" and the old
# row 22 "质量" ends up at row 21).
$ws.Range("A7").Value = "PROTEIN"
$ws.Range("A9").Value = "ORGANIZATION (ENVIRONMENT)"

$ws.Range("A11").Value = "MATTEr"
$ws.Range("B11").Value = $ws.Range("B12").Value2
$ws.Range("C11").Value = $ws.Range("C12").Value2

$ws.Range("A12").Value = "LOC"
$ws.Range("A13").Value = "MAT"
$ws.Range("A14").Value = "FIELD"
$ws.Range("A15").Value = "DATE"
$ws.Range("A16").Value = "COLUMN"
$ws.Range("A17").Value = "FUNCTION"
$ws.Range("A18").Value = "REALQUALITYTY"
$ws.Range("A19").Value = "REGULATIONS"
$ws.Range("A20").Value = "STRING"
$ws.Range("A21").Value = "质量"

# Remove the now-duplicate last row (22), shifting cells up.
$ws.Rows.Item(22).Delete()
